$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.057.08'
$ws.Range("E2").Value = '  -5.40%  '
$ws.Range("D3").Value = '2.174.37'
$ws.Range("E3").Value = '  -8.02%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''477.34'
$ws.Range("E5").Value = '  -4.65%  '
$ws.Range("D6").Value = '''122.46'
$ws.Range("E6").Value = '  -4.89%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '''0.511'
$ws.Range("E8").Value = '  -6.02%  '
$ws.Range("D9").Value = '2.182.92'
$ws.Range("E9").Value = '  -7.79%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.0902'
$ws.Range("E10").Value = '  -7.91%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.147'
$ws.Range("E11").Value = '  -2.19%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.308'
$ws.Range("E12").Value = '  -4.33%  '
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").Value = '''4.51'
$ws.Range("E13").Value = '  -7.18%  '
$ws.Range("D14").Value = '2.562.31'
$ws.Range("E14").Value = '  -7.96%  '
$ws.Range("D15").Value = '''20.72'
$ws.Range("E15").Value = '  -3.14%  '
$ws.Range("D16").Value = '52.978.11'
$ws.Range("E16").Value = '  -5.49%  '
$ws.Range("D17").Value = '''0.0000125'
$ws.Range("E17").Value = '  -5.11%  '
$ws.Range("D18").Value = '2.171.63'
$ws.Range("E18").Value = '  -7.55%  '
$ws.Range("D19").Value = '''9.42'
$ws.Range("E19").Value = '  -5.66%  '
$ws.Range("D20").Value = '''3.88'
$ws.Range("E20").Value = '  -3.92%  '
$ws.Range("D21").Value = '''290.59'
$ws.Range("E21").Value = '  -5.22%  '
$ws.Range("D22").Value = '''0.998'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '''5.96'
$ws.Range("E23").Value = '  -5.07%  '
$ws.Range("D24").Value = '''61.96'
$ws.Range("E24").Value = '  -6.05%  '
$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").Value = '''0.359'
$ws.Range("E26").Value = '  -2.43%  '
$ws.Range("D27").Value = '2.277.51'
$ws.Range("E27").Value = '  -7.98%  '
$ws.Range("D28").Value = '''0.142'
$ws.Range("E28").Value = '  -3.07%  '
$ws.Range("D29").Value = '''6.88'
$ws.Range("E29").Value = '  -4.36%  '
$ws.Range("D30").Value = '''165.82'
$ws.Range("E30").Value = '  -3.12%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  -5.15%  '
$ws.Range("D33").Value = '''0.996'
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("D34").Value = '0.0₃0644'
$ws.Range("E34").Value = '  -9.02%  '
$ws.Range("D35").Value = '''5.55'
$ws.Range("E35").Value = '  -3.02%  '
$ws.Range("D36").Value = '''1.02'
$ws.Range("E36").Value = '  -5.23%  '
$ws.Range("D37").Value = '''16.98'
$ws.Range("E37").Value = '  -3.47%  '
$ws.Range("E38").Value = '  -4.05%  '
$ws.Range("D39").Value = '''0.799'
$ws.Range("E39").Value = '  +1.30%  '
$ws.Range("D40").Value = '''35.36'
$ws.Range("D41").Value = '''3.48'
$ws.Range("E41").Value = '  -6.60%  '
$ws.Range("D42").Value = '''0.361'
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("E43").Value = '  -3.54%  '
$ws.Range("D44").Value = '''3.19'
$ws.Range("E44").Value = '  -4.51%  '
$ws.Range("D45").Value = '''120.90'
$ws.Range("E45").Value = '  -6.33%  '
$ws.Range("D46").Value = '''4.61'
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").Value = '''0.0864'
$ws.Range("E47").Value = '  -3.94%  '
$ws.Range("D48").Value = '''0.523'
$ws.Range("E48").Value = '  -7.21%  '
$ws.Range("D49").Value = '''0.0462'
$ws.Range("E49").Value = '  -3.72%  '
$ws.Range("D50").Value = '''224.74'
$ws.Range("E50").Value = '  -6.01%  '
$ws.Range("D51").Value = '''0.0197'
$ws.Range("E51").Value = '  -4.67%  '
